# OLX Monitor 2026-02-19 08:46 — append newly discovered listings
# to the detail log sheet (rows 75-82), matching the existing
# "timestamp / profile / title / price / added-date / days / url / slug"
# row layout (and styles) already used through row 74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 8 new rows, in column order A..H.
$rows = @(
    @("2026-02-19 08:46:44", "poqui", "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy", 2499, "28.10.2025", 113, "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html", "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"),
    @("2026-02-19 08:46:44", "poqui", "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda", 2299, "19.01.2026", 30, "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html", "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"),
    @("2026-02-19 08:46:44", "poqui", "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza", 2049, "19.12.2025", 61, "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html", "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"),
    @("2026-02-19 08:46:44", "poqui", "Przytulny pokój blisko Politechniki – ul. Przytulna", 599, "10.10.2025", 131, "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html", "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"),
    @("2026-02-19 08:46:44", "pokojewlublinie", "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58", 58640, "11.08.2025", 191, "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html", "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"),
    @("2026-02-19 08:46:44", "pokojewlublinie", "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12", 12640, "19.01.2026", 30, "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html", "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"),
    @("2026-02-19 08:46:44", "dawnypatron", "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.", 730, "20.09.2024", 516, "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html", "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"),
    @("2026-02-19 08:46:44", "dawnypatron", "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14", 14690, "05.12.2025", 75, "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html", "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv")
)

# The "Dni" (days) column F uses a highlight style (s=15, red font) for
# most rows, but a plain centered style (s=14) when the value is the
# "just refreshed" 30-day marker. Mirrors the existing rows' pattern.
$fStyleRow = @(74, 7, 74, 74, 74, 7, 74, 74)

# Scratch cell used to push "DD.MM.YYYY" strings in as literal text.
# Excel's smart-input would otherwise reinterpret day-ambiguous dates
# (day <= 12, so also valid as MM.DD.YYYY) as real date serials; forcing
# the scratch cell to Text format first and then copying only the
# *value* across keeps the destination cell's already-applied style
# (from the row clone below) untouched. NumberFormat is re-applied
# before every use since Clear() resets it along with the old value.
$scratch = $ws.Range("Z1000")

$startRow = 75
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Clone the formatting of the last existing data row onto the new
    # row first, so alignment/styles match exactly, then overwrite F's
    # style source per-row to match style 14 vs 15 as needed.
    $ws.Range("A74:H74").Copy()
    $ws.Range("A$r`:H$r").PasteSpecial(-4122)

    $fSrc = $fStyleRow[$i]
    if ($fSrc -ne 74) {
        $ws.Range("F$fSrc").Copy()
        $ws.Range("F$r").PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]

    $day = [int]($data[4].Substring(0, 2))
    if ($day -le 12) {
        # Ambiguous as MM.DD.YYYY too -> force literal text via scratch cell.
        $scratch.NumberFormat = "@"
        $scratch.Value = $data[4]
        $scratch.Copy()
        $ws.Cells.Item($r, 5).PasteSpecial(-4163)
        $scratch.Clear()
    } else {
        $ws.Cells.Item($r, 5).Value = $data[4]
    }

    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
